$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Locate the target paragraph (the one beginning with "From our
# preliminary studies...") so later logic doesn't depend on a fixed
# paragraph index.
# -----------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs.Item($i).Range.Text
    if ($ptxt -like "From our preliminary studies*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph"
}

# -----------------------------------------------------------------
# 1) "...presence of higher order interaction " -> insert "(HOI) "
#    right after it (before "help").
# -----------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("higher order interaction ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPoint = $d.Range($rng.End, $rng.End)
$insPoint.InsertAfter("(HOI) ")

# -----------------------------------------------------------------
# 2) "help" -> "helps" (insert "s" right after "help").
# -----------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("help", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPoint = $d.Range($rng.End, $rng.End)
$insPoint.InsertAfter("s")

# -----------------------------------------------------------------
# 3) Append "This forms our first hypothesis." to the end of the
#    paragraph (after "...in the community. ").
# -----------------------------------------------------------------
$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range
$endPoint = $d.Range($r.End - 1, $r.End - 1)
$endPoint.InsertAfter("This forms our first hypothesis.")

# -----------------------------------------------------------------
# 4) Insert a new, empty paragraph right after this paragraph
#    (before the pre-existing empty paragraph that followed it).
# -----------------------------------------------------------------
$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range
$r.InsertParagraphAfter()
